# Worker List.xlsx - "Add files via upload" commit
# Adds a new "S/N" (serial number) column D, with values for a handful of
# workers, plus matching "ID" values for those same rows that previously had
# no ID recorded. Also tidies up the sheet view / page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New column D: "S/N"
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 16.75

$ws.Range("D1").Value = "S/N"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Row 2 - Adam Price
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "86FA1A07"
$ws.Range("C2").HorizontalAlignment = -4131
$ws.Range("D2").Value = "N521D5060024"

# ---------------------------------------------------------------------
# Row 3 - Amit Mallik
# ---------------------------------------------------------------------
$ws.Range("C3").Value = "8166BF3C"
$ws.Range("C3").HorizontalAlignment = -4131
$ws.Range("D3").Value = "N521D5060019"

# ---------------------------------------------------------------------
# Row 9 - Damarley Barrett
# ---------------------------------------------------------------------
$ws.Range("C9").Value = "B93BAE45"
$ws.Range("C9").HorizontalAlignment = -4131
$ws.Range("D9").Value = "N521D5060035"

# ---------------------------------------------------------------------
# Row 10 - Dominic Miller
# ---------------------------------------------------------------------
$ws.Range("C10").Value = "DBFC40AA"
$ws.Range("D10").Value = "N521D5060033"

# ---------------------------------------------------------------------
# Row 34 - Sayed Hossiny
# ---------------------------------------------------------------------
$ws.Range("C34").Value = "BCAFEC01"
$ws.Range("D34").Value = "N521D5060027"

# ---------------------------------------------------------------------
# View / page tidy-up
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$null = $ws.Range("G30").Select()

Write-Host "Worker List updated: added S/N column and new ID/S-N values."
